# clientes.xlsx: drop the "tipo" and "anotaciones" columns, uppercase the
# client name, and reset "monto" back to 0 for the first record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column G ("anotaciones") first so column letters to its left
# (in particular C, "tipo") are unaffected by this deletion.
$ws.Columns("G").Delete()

# Remove column C ("tipo"); telefono/direccion/monto shift left into
# C/D/E.
$ws.Columns("C").Delete()

# Data fixes on the remaining row 2.
$ws.Range("B2").Value = "FRANCO BONIN"
$ws.Range("E2").Value = 0
